$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws.Range("D61").Value = 5
Write-Host "done"
